# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.480.26"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "'3.445.51"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'584.10"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "'176.49"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "'3.438.51"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  -5.66%  "
$ws.Range("D11").Value = "'6.87"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "'0.417"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").Value = "'4.043.62"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "'30.24"
$ws.Range("E14").Value = "  -4.74%  "
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "'66.472.96"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "'3.453.32"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'5.95"
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("D20").Value = "'13.82"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'379.13"
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").Value = "'7.81"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.538"
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'72.25"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "'5.74"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'24.31"
$ws.Range("E31").Value = "  +3.82%  "
$ws.Range("E32").Value = "  -4.87%  "
$ws.Range("D33").Value = "'1.99"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -5.71%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'7.18"
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "'160.87"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "'29.37"
$ws.Range("E39").Value = "  +12.41%  "
$ws.Range("D40").Value = "'0.890"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.64"
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D43").Value = "'4.51"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'2.733.06"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "'6.39"
$ws.Range("E45").Value = "  -5.64%  "
$ws.Range("D46").Value = "'0.0696"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "'40.56"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'24.59"
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("D49").Value = "'0.0293"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'308.80"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "'0.829"
$ws.Range("E51").Value = "  -1.31%  "
